# ==================================================================
# feat: add 2022-Q4 data
#
# 1) "总计" (Total) sheet: push the existing rows down and add a new
#    leading row for "2022-Q4".
# 2) Insert a brand-new "2022-Q4" worksheet (with per-fund detail)
#    right before the existing "2022-Q3" worksheet.
# ==================================================================

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item("总计")

# --- 1) "总计" sheet: shift rows down, insert 2022-Q4 at the top ---
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 23
$total.Range("D2").Value = 1.01

$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 18
$total.Range("D3").Value = 0.92

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.29

# Carry the "index column" look (bold + border + centred) onto the new row
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

# --- 2) Insert the new "2022-Q4" worksheet, right before "2022-Q3" ---
$existingQ3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($existingQ3)
$q4.Name = "2022-Q4"
# Re-fetch "2022-Q3" now that the sheet collection has changed, so the
# style-source range below references a live object.
$q3 = $wb.Worksheets.Item("2022-Q3")

# Header row (row 1), columns B:H
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data rows 2-24
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'006234"
$q4.Range("C2").Value = "万家汽车新趋势混合C"
$q4.Range("D2").Value = "'4.06"
$q4.Range("E2").Value = "'90.27"
$q4.Range("F2").Value = "'3.95"
$q4.Range("G2").Value = "'0.1604"
$q4.Range("H2").Value = 3

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'003516"
$q4.Range("C3").Value = "国泰融安多策略灵活配置混合A"
$q4.Range("D3").Value = "'7.77"
$q4.Range("E3").Value = "'78.85"
$q4.Range("F3").Value = "'2.02"
$q4.Range("G3").Value = "'0.1570"
$q4.Range("H3").Value = 7

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'006132"
$q4.Range("C4").Value = "万家智造优势混合A"
$q4.Range("D4").Value = "'4.10"
$q4.Range("E4").Value = "'93.92"
$q4.Range("F4").Value = "'3.24"
$q4.Range("G4").Value = "'0.1328"
$q4.Range("H4").Value = 9

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'010335"
$q4.Range("C5").Value = "华宝竞争优势混合"
$q4.Range("D5").Value = "'2.59"
$q4.Range("E5").Value = "'94.54"
$q4.Range("F5").Value = "'4.95"
$q4.Range("G5").Value = "'0.1282"
$q4.Range("H5").Value = 6

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'005983"
$q4.Range("C6").Value = "上投摩根核心精选股票A"
$q4.Range("D6").Value = "'3.39"
$q4.Range("E6").Value = "'94.31"
$q4.Range("F6").Value = "'3.28"
$q4.Range("G6").Value = "'0.1112"
$q4.Range("H6").Value = 10

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'020023"
$q4.Range("C7").Value = "国泰事件驱动策略混合A"
$q4.Range("D7").Value = "'2.19"
$q4.Range("E7").Value = "'82.03"
$q4.Range("F7").Value = "'3.34"
$q4.Range("G7").Value = "'0.0731"
$q4.Range("H7").Value = 5

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "'006233"
$q4.Range("C8").Value = "万家汽车新趋势混合A"
$q4.Range("D8").Value = "'1.65"
$q4.Range("E8").Value = "'90.27"
$q4.Range("F8").Value = "'3.95"
$q4.Range("G8").Value = "'0.0652"
$q4.Range("H8").Value = 3

$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "'012880"
$q4.Range("C9").Value = "国泰景气优选混合A"
$q4.Range("D9").Value = "'3.24"
$q4.Range("E9").Value = "'75.90"
$q4.Range("F9").Value = "'1.99"
$q4.Range("G9").Value = "'0.0645"
$q4.Range("H9").Value = 8

$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "'009379"
$q4.Range("C10").Value = "中银成长优选股票A"
$q4.Range("D10").Value = "'1.51"
$q4.Range("E10").Value = "'83.26"
$q4.Range("F10").Value = "'2.18"
$q4.Range("G10").Value = "'0.0329"
$q4.Range("H10").Value = 8

$q4.Range("A11").Value = 9
$q4.Range("B11").Value = "'006133"
$q4.Range("C11").Value = "万家智造优势混合C"
$q4.Range("D11").Value = "'0.78"
$q4.Range("E11").Value = "'93.92"
$q4.Range("F11").Value = "'3.24"
$q4.Range("G11").Value = "'0.0253"
$q4.Range("H11").Value = 9

$q4.Range("A12").Value = 10
$q4.Range("B12").Value = "'005728"
$q4.Range("C12").Value = "华宝绿色主题混合A"
$q4.Range("D12").Value = "'0.40"
$q4.Range("E12").Value = "'91.97"
$q4.Range("F12").Value = "'4.73"
$q4.Range("G12").Value = "'0.0189"
$q4.Range("H12").Value = 6

$q4.Range("A13").Value = 11
$q4.Range("B13").Value = "'001231"
$q4.Range("C13").Value = "银华泰利灵活配置混合A"
$q4.Range("D13").Value = "'0.87"
$q4.Range("E13").Value = "'22.82"
$q4.Range("F13").Value = "'0.95"
$q4.Range("G13").Value = "'0.0083"
$q4.Range("H13").Value = 5

$q4.Range("A14").Value = 12
$q4.Range("B14").Value = "'005437"
$q4.Range("C14").Value = "易方达易百智能量化策略灵活配置混合A"
$q4.Range("D14").Value = "'0.76"
$q4.Range("E14").Value = "'94.68"
$q4.Range("F14").Value = "'0.99"
$q4.Range("G14").Value = "'0.0075"
$q4.Range("H14").Value = 10

$q4.Range("A15").Value = 13
$q4.Range("B15").Value = "'005000"
$q4.Range("C15").Value = "泰康泉林量化价值精选混合A"
$q4.Range("D15").Value = "'0.31"
$q4.Range("E15").Value = "'89.21"
$q4.Range("F15").Value = "'1.65"
$q4.Range("G15").Value = "'0.0051"
$q4.Range("H15").Value = 7

$q4.Range("A16").Value = 14
$q4.Range("B16").Value = "'012881"
$q4.Range("C16").Value = "国泰景气优选混合C"
$q4.Range("D16").Value = "'0.22"
$q4.Range("E16").Value = "'75.90"
$q4.Range("F16").Value = "'1.99"
$q4.Range("G16").Value = "'0.0044"
$q4.Range("H16").Value = 8

$q4.Range("A17").Value = 15
$q4.Range("B17").Value = "'014960"
$q4.Range("C17").Value = "国泰融安多策略灵活配置混合C"
$q4.Range("D17").Value = "'0.14"
$q4.Range("E17").Value = "'78.85"
$q4.Range("F17").Value = "'2.02"
$q4.Range("G17").Value = "'0.0028"
$q4.Range("H17").Value = 7

$q4.Range("A18").Value = 16
$q4.Range("B18").Value = "'014455"
$q4.Range("C18").Value = "中银成长优选股票C"
$q4.Range("D18").Value = "'0.12"
$q4.Range("E18").Value = "'83.26"
$q4.Range("F18").Value = "'2.18"
$q4.Range("G18").Value = "'0.0026"
$q4.Range("H18").Value = 8

$q4.Range("A19").Value = 17
$q4.Range("B19").Value = "'005111"
$q4.Range("C19").Value = "泰康泉林量化价值精选混合C"
$q4.Range("D19").Value = "'0.14"
$q4.Range("E19").Value = "'89.21"
$q4.Range("F19").Value = "'1.65"
$q4.Range("G19").Value = "'0.0023"
$q4.Range("H19").Value = 7

$q4.Range("A20").Value = 18
$q4.Range("B20").Value = "'012799"
$q4.Range("C20").Value = "华宝绿色主题混合C"
$q4.Range("D20").Value = "'0.04"
$q4.Range("E20").Value = "'91.97"
$q4.Range("F20").Value = "'4.73"
$q4.Range("G20").Value = "'0.0019"
$q4.Range("H20").Value = 6

$q4.Range("A21").Value = 19
$q4.Range("B21").Value = "'014937"
$q4.Range("C21").Value = "上投摩根核心精选股票C"
$q4.Range("D21").Value = "'0.05"
$q4.Range("E21").Value = "'94.31"
$q4.Range("F21").Value = "'3.28"
$q4.Range("G21").Value = "'0.0016"
$q4.Range("H21").Value = 10

$q4.Range("A22").Value = 20
$q4.Range("B22").Value = "'005438"
$q4.Range("C22").Value = "易方达易百智能量化策略灵活配置混合C"
$q4.Range("D22").Value = "'0.15"
$q4.Range("E22").Value = "'94.68"
$q4.Range("F22").Value = "'0.99"
$q4.Range("G22").Value = "'0.0015"
$q4.Range("H22").Value = 10

$q4.Range("A23").Value = 21
$q4.Range("B23").Value = "'002328"
$q4.Range("C23").Value = "银华泰利灵活配置混合C"
$q4.Range("D23").Value = "'0.03"
$q4.Range("E23").Value = "'22.82"
$q4.Range("F23").Value = "'0.95"
$q4.Range("G23").Value = "'0.0003"
$q4.Range("H23").Value = 5

$q4.Range("A24").Value = 22
$q4.Range("B24").Value = "'015592"
$q4.Range("C24").Value = "国泰事件驱动策略混合C"
$q4.Range("D24").Value = "'0.01"
$q4.Range("E24").Value = "'82.03"
$q4.Range("F24").Value = "'3.34"
$q4.Range("G24").Value = "'0.0003"
$q4.Range("H24").Value = 5

# Match the formatting used by the other quarterly sheets:
#   - header row (B1:H1) and index column (A2:A24) are bold, bordered,
#     center-aligned, taken from the "2022-Q3" sheet which already has it
$cols = @("B","C","D","E","F","G","H")
foreach ($col in $cols) {
    $q3.Range("$col" + "1").Copy()
    $q4.Range("$col" + "1").PasteSpecial(-4122)
}
for ($r = 2; $r -le 24; $r++) {
    $q3.Range("A2").Copy()
    $q4.Range("A$r").PasteSpecial(-4122)
}
